$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on every Price/Volume cell being updated so Excel
# does not auto-convert numeric-looking strings (e.g. "0.111", "5.76")
# into floating point numbers -- these columns are plain text in the source data.
$textCells = @("D2", "E2", "D3", "E3", "D5", "E5", "D6", "E6", "D8", "E8", "E9", "E10", "D11", "E11", "E12", "E13", "D14", "E14", "D15", "E15", "E16", "D17", "E17", "D18", "E18", "E19", "D20", "E20", "E21", "E22", "E23", "E24", "D25", "E25", "E26", "E27", "D28", "E28", "E29", "E30", "D31", "E31", "D32", "E32", "D33", "E33", "D34", "E34", "E35", "E36", "D37", "E37", "E38", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "E45", "E46", "D47", "E47", "D48", "E48", "E50", "E51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "63.786.18"
$ws.Range("E2").Value = "  +0.65%  "
$ws.Range("D3").Value = "3.144.30"
$ws.Range("E3").Value = "  +0.91%  "
$ws.Range("D5").Value = "586.69"
$ws.Range("E5").Value = "  +0.18%  "
$ws.Range("D6").Value = "145.49"
$ws.Range("E6").Value = "  +0.24%  "
$ws.Range("D8").Value = "3.135.61"
$ws.Range("E8").Value = "  +0.89%  "
$ws.Range("E9").Value = "  -0.24%  "
$ws.Range("E10").Value = "  +6.47%  "
$ws.Range("D11").Value = "5.76"
$ws.Range("E11").Value = "  -0.50%  "
$ws.Range("E12").Value = "  -2.20%  "
$ws.Range("E13").Value = "  -0.37%  "
$ws.Range("D14").Value = "36.94"
$ws.Range("E14").Value = "  +3.52%  "
$ws.Range("D15").Value = "3.665.48"
$ws.Range("E15").Value = "  +1.00%  "
$ws.Range("E16").Value = "  -1.38%  "
$ws.Range("D17").Value = "3.142.05"
$ws.Range("E17").Value = "  +0.93%  "
$ws.Range("D18").Value = "63.565.26"
$ws.Range("E18").Value = "  +0.47%  "
$ws.Range("E19").Value = "  -1.17%  "
$ws.Range("D20").Value = "463.52"
$ws.Range("E20").Value = "  -0.95%  "
$ws.Range("E21").Value = "  +0.64%  "
$ws.Range("E22").Value = "  +0.22%  "
$ws.Range("E23").Value = "  -1.51%  "
$ws.Range("E24").Value = "  -2.58%  "
$ws.Range("D25").Value = "81.32"
$ws.Range("E25").Value = "  -0.97%  "
$ws.Range("E26").Value = "  +0.60%  "
$ws.Range("E27").Value = "  -0.31%  "
$ws.Range("D28").Value = "9.09"
$ws.Range("E28").Value = "  +6.39%  "
$ws.Range("E29").Value = "  -0.19%  "
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("D31").Value = "2.21"
$ws.Range("E31").Value = "  -1.14%  "
$ws.Range("D32").Value = "6.97"
$ws.Range("E32").Value = "  +1.20%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "0.111"
$ws.Range("E33").Value = "  +0.69%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "26.92"
$ws.Range("E34").Value = "  -0.24%  "
$ws.Range("E35").Value = "  -1.19%  "
$ws.Range("E36").Value = "  -0.76%  "
$ws.Range("D37").Value = "3.38"
$ws.Range("E37").Value = "  +1.23%  "
$ws.Range("E38").Value = "  -4.98%  "
$ws.Range("E39").Value = "  -1.37%  "
$ws.Range("D40").Value = "50.68"
$ws.Range("E40").Value = "  +0.18%  "
$ws.Range("D41").Value = "440.54"
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("D42").Value = "8.76"
$ws.Range("E42").Value = "  +0.26%  "
$ws.Range("D43").Value = "0.0370"
$ws.Range("E43").Value = "  +0.24%  "
$ws.Range("D44").Value = "2.908.58"
$ws.Range("E44").Value = "  -0.69%  "
$ws.Range("E45").Value = "  -2.23%  "
$ws.Range("E46").Value = "  -2.28%  "
$ws.Range("D47").Value = "36.45"
$ws.Range("E47").Value = "  +2.45%  "
$ws.Range("D48").Value = "125.68"
$ws.Range("E48").Value = "  +1.97%  "
$ws.Range("E50").Value = "  -0.77%  "
$ws.Range("E51").Value = "  -1.65%  "
